# Applies the "Corrected excel sheets for application fix issues" edit:
#   - Summary sheet: shrink the lingering selection by one row.
#   - Repayment schedule sheet: add the missing "O" column values
#     (rows 2-15) mirroring the existing N/P columns, style included.
#   - Transactions sheet: renumber the first three transaction IDs and
#     move the active selection.

$wb = $excel.ActiveWorkbook

# ---- Summary sheet -------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()
$wsSummary.Range("A7:XFD14").Select()

# ---- Repayment schedule sheet ---------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Data rows (3,5,6,7,8,9,10,11,12,13,15) get a numeric 0 in column O,
# matching the existing N/P columns for that row. Rows 2,4,14 only had
# an (empty) style-only cell in N/P, so column O mirrors that: style
# applied, no value.
$oZeroRows = @(3,5,6,7,8,9,10,11,12,13,15)
$oEmptyRows = @(2,4,14)

$xlPasteFormats = -4122

foreach ($r in $oZeroRows) {
    $wsSchedule.Cells.Item($r, 14).Copy()
    $wsSchedule.Cells.Item($r, 15).PasteSpecial($xlPasteFormats)
    $wsSchedule.Cells.Item($r, 15).Value = 0
}

foreach ($r in $oEmptyRows) {
    $wsSchedule.Cells.Item($r, 14).Copy()
    $wsSchedule.Cells.Item($r, 15).PasteSpecial($xlPasteFormats)
}

# ---- Transactions sheet ----------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")

$wsTransactions.Cells.Item(2, 1).Value = 68
$wsTransactions.Cells.Item(3, 1).Value = 67
$wsTransactions.Cells.Item(4, 1).Value = 66

$wsTransactions.Activate()
$wsTransactions.Range("A2:L4").Select()
